$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1043256666666667
$ws.Range("H2").Value = 0.312977
$ws.Range("I2").Value = 0.02547563162231953
$ws.Range("J2").Value = 0.02547563162231953
$ws.Range("M2").Value = 14.440165
$ws.Range("N2").Value = 43.320495
$ws.Range("O2").Value = 0.1441015470002482
$ws.Range("P2").Value = 0.1441015470002482
$ws.Range("Q2").Value = 1.506479840401667
$ws.Range("R2").Value = 13.558318563615
$ws.Range("S2").Value = 0.003671077927584688
$ws.Range("T2").Value = 0.003671077927584687
# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1043256666666667
$ws.Range("H3").Value = 0.312977
$ws.Range("I3").Value = 0.02547563162231953
$ws.Range("J3").Value = 0.02547563162231953
$ws.Range("O3").Value = 0.3846359116098663
$ws.Range("P3").Value = 0.3846359116098662
$ws.Range("Q3").Value = 4.021096641896444
$ws.Range("R3").Value = 36.189869777068
$ws.Range("S3").Value = 0.009798842792888009
$ws.Range("T3").Value = 0.009798842792888007
# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1043256666666667
$ws.Range("H4").Value = 0.312977
$ws.Range("I4").Value = 0.02547563162231953
$ws.Range("J4").Value = 0.02547563162231953
$ws.Range("M4").Value = 21.954262
$ws.Range("N4").Value = 65.862786
$ws.Range("O4").Value = 0.2190863551385157
$ws.Range("P4").Value = 0.2190863551385156
$ws.Range("Q4").Value = 2.290393019324667
$ws.Range("R4").Value = 20.613537173922
$ws.Range("S4").Value = 0.005581363276985497
$ws.Range("T4").Value = 0.005581363276985496
# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1043256666666667
$ws.Range("H5").Value = 0.312977
$ws.Range("I5").Value = 0.02547563162231953
$ws.Range("J5").Value = 0.02547563162231953
$ws.Range("M5").Value = 25.27013633333333
$ws.Range("N5").Value = 75.81040899999999
$ws.Range("O5").Value = 0.2521761862513699
$ws.Range("P5").Value = 0.2521761862513699
$ws.Range("Q5").Value = 2.636323819732555
$ws.Range("R5").Value = 23.726914377593
$ws.Range("S5").Value = 0.00642434762486134
$ws.Range("T5").Value = 0.00642434762486134
# Row 6
$ws.Range("I6").Value = 0.9745243683776804
$ws.Range("J6").Value = 0.9745243683776804
$ws.Range("M6").Value = 14.440165
$ws.Range("N6").Value = 43.320495
$ws.Range("O6").Value = 0.1441015470002482
$ws.Range("P6").Value = 0.1441015470002482
$ws.Range("Q6").Value = 57.62767089373833
$ws.Range("R6").Value = 518.649038043645
$ws.Range("S6").Value = 0.1404304690726635
$ws.Range("T6").Value = 0.1404304690726635
# Row 7
$ws.Range("I7").Value = 0.9745243683776804
$ws.Range("J7").Value = 0.9745243683776804
$ws.Range("O7").Value = 0.3846359116098663
$ws.Range("P7").Value = 0.3846359116098662
$ws.Range("S7").Value = 0.3748370688169783
$ws.Range("T7").Value = 0.3748370688169782
# Row 8
$ws.Range("I8").Value = 0.9745243683776804
$ws.Range("J8").Value = 0.9745243683776804
$ws.Range("M8").Value = 21.954262
$ws.Range("N8").Value = 65.862786
$ws.Range("O8").Value = 0.2190863551385157
$ws.Range("P8").Value = 0.2190863551385156
$ws.Range("Q8").Value = 87.61485656506733
$ws.Range("R8").Value = 788.5337090856059
$ws.Range("S8").Value = 0.2135049918615302
$ws.Range("T8").Value = 0.2135049918615301
# Row 9
$ws.Range("I9").Value = 0.9745243683776804
$ws.Range("J9").Value = 0.9745243683776804
$ws.Range("M9").Value = 25.27013633333333
$ws.Range("N9").Value = 75.81040899999999
$ws.Range("O9").Value = 0.2521761862513699
$ws.Range("P9").Value = 0.2521761862513699
$ws.Range("Q9").Value = 100.8478158010821
$ws.Range("R9").Value = 907.6303422097388
$ws.Range("S9").Value = 0.2457518386265086
$ws.Range("T9").Value = 0.2457518386265086
